# "Misc change in Consignment Entry"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "isEmpty" column (C) used to hold TRUE/FALSE booleans; switch it to
# the more human-readable "Y"/"N" text values.
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"

# Resize the saved workbook window a bit taller.
$excel.ActiveWindow.Height = 7650

# Move the active selection to C10.
$ws.Range("C10").Select() | Out-Null
